$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers need to be forced
# to Text format first, otherwise Excel auto-converts the string to a number
# (these columns are textual price strings, not numeric values).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "34.148.99"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").Value = "1.778.66"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "225.85"
$ws.Range("E5").Value = "  +0.59%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "31.64"
$ws.Range("E8").Value = "  -0.34%  "

$ws.Range("E9").Value = "  +0.98%  "

$ws.Range("D10").Value = "0.0692"
$ws.Range("E10").Value = "  +2.00%  "

$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("D12").Value = "2.034.37"
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").Value = "1.788.98"
$ws.Range("E13").Value = "  -0.63%  "

$ws.Range("E14").Value = "  -1.70%  "

$ws.Range("D15").Value = "34.109.68"

$ws.Range("D16").Value = "0.622"
$ws.Range("E16").Value = "  +2.00%  "

$ws.Range("D17").Value = "4.19"
$ws.Range("E17").Value = "  +1.22%  "

$ws.Range("D18").Value = "67.85"

$ws.Range("D19").Value = "0.0₃0799"
$ws.Range("E19").Value = "  +3.73%  "

$ws.Range("D20").Value = "245.37"
$ws.Range("E20").Value = "  +2.59%  "

$ws.Range("D21").Value = "10.97"
$ws.Range("E21").Value = "  +4.03%  "

$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("E23").Value = "  +1.75%  "

$ws.Range("E24").Value = "  -0.99%  "

$ws.Range("D25").Value = "162.08"
$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("D26").Value = "7.20"
$ws.Range("E26").Value = "  +2.43%  "

$ws.Range("D27").Value = "16.29"
$ws.Range("E27").Value = "  +1.24%  "

$ws.Range("E28").Value = "  +1.81%  "

$ws.Range("E29").Value = "  +0.33%  "

$ws.Range("E30").Value = "  +0.71%  "

$ws.Range("D31").Value = "0.0520"
$ws.Range("E31").Value = "  +2.18%  "

$ws.Range("E32").Value = "  +3.92%  "

$ws.Range("D33").Value = "3.71"
$ws.Range("E33").Value = "  +5.34%  "

$ws.Range("E34").Value = "  -1.37%  "

$ws.Range("D35").Value = "1.438.60"
$ws.Range("E35").Value = "  +3.47%  "

$ws.Range("E36").Value = "  +3.87%  "

$ws.Range("E37").Value = "  +6.40%  "

$ws.Range("E38").Value = "  +2.44%  "

$ws.Range("E39").Value = "  +0.33%  "

$ws.Range("D40").Value = "80.15"
$ws.Range("E40").Value = "  +2.37%  "

$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("D42").Value = "0.921"
$ws.Range("E42").Value = "  +0.87%  "

$ws.Range("E43").Value = "  +0.69%  "

$ws.Range("D44").Value = "13.46"
$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("E46").Value = "  +3.61%  "

$ws.Range("E47").Value = "  +0.14%  "

$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("D49").Value = "1.937.00"

$ws.Range("D50").Value = "104.15"
$ws.Range("E50").Value = "  -1.12%  "

$ws.Range("E51").Value = "  +0.15%  "
